# Fix soft assertion issue: add an "Index" column to the FooterLinks sheet
# so each row is uniquely identifiable, shifting the existing Title /
# Link Text / Link Route columns one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Insert a new blank column before column A; this shifts the existing
# A/B/C columns (Title/Link Text/Link Route) into B/C/D and preserves
# their formatting, shared-string values and per-cell styles.
$ws.Columns.Item(1).Insert()

# Header for the new index column.
$ws.Range("A1").Value = "Index"

# Numeric index values for each data row (section.item numbering).
$indexValues = @(1.1, 1.2, 1.3, 1.4, 1.5, 1.6, 2.1, 2.2, 2.3, 2.4, 2.5, 3.1, 3.2, 3.3, 3.4, 3.5, 3.6)

for ($i = 0; $i -lt $indexValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $indexValues[$i]
}

# Restore selection to match the saved view state.
$ws.Range("D16").Select()
